# Generate Report for Handback
# Populates the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" / "Error Detail" columns (I, J, K, P) for the
# 35ca43d0-4c2a-49d3-8c21-9052e9f7362e row (row 5) on both the zh-cn and
# de-de worksheets, now that a (stale) handback has come in for that file.

$wb = $excel.ActiveWorkbook

$targetMdDisplay = "35ca43d0-4c2a-49d3-8c21-9052e9f7362e.md"
$targetMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68df7fb50e16849dc1387e560e998ef721dec850/e2e/35ca43d0-4c2a-49d3-8c21-9052e9f7362e.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6568e4e0c75a0dd3837e59c4126497d1daf82e5e/e2e/35ca43d0-4c2a-49d3-8c21-9052e9f7362e.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/68df7fb50e16849dc1387e560e998ef721dec850/e2e/35ca43d0-4c2a-49d3-8c21-9052e9f7362e.md."

# Blue colour (FF6495ED, i.e. RGB 100,149,237) used by the workbook's
# existing "HyperLink" style, expressed as an OLE BGR long for Font.Color.
$hyperlinkColor = 15570276

# Width of 40 "characters" after Excel's internal +~0.83 padding -- this is
# the same stored width ("40") already used by several other columns in
# this workbook.
$wideColumnWidth = 39.1666666666667

$zhCnTargetFile = "35ca43d0-4c2a-49d3-8c21-9052e9f7362e.133f7d495c2a2758699fa96f6c1df0ff55e565de.zh-cn.xlf"
$zhCnHandbackDate = "2016-10-19 11:15:50"

$deDeTargetFile = "35ca43d0-4c2a-49d3-8c21-9052e9f7362e.133f7d495c2a2758699fa96f6c1df0ff55e565de.de-de.xlf"
$deDeHandbackDate = "2016-10-19 11:16:09"

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$iCellZhCn = $wsZhCn.Range("I5")
$iCellZhCn.Value = $targetMdDisplay
$wsZhCn.Hyperlinks.Add($iCellZhCn, $targetMdUrl, [Type]::Missing, [Type]::Missing, $targetMdDisplay) | Out-Null
$iCellZhCn.Font.Underline = 2
$iCellZhCn.Font.Color = $hyperlinkColor

$wsZhCn.Range("J5").Value = $zhCnTargetFile
$wsZhCn.Range("K5").Value = $zhCnHandbackDate
$wsZhCn.Range("P5").Value = $errorDetail

$wsZhCn.Columns.Item(16).ColumnWidth = $wideColumnWidth

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")

$iCellDeDe = $wsDeDe.Range("I5")
$iCellDeDe.Value = $targetMdDisplay
$wsDeDe.Hyperlinks.Add($iCellDeDe, $targetMdUrl, [Type]::Missing, [Type]::Missing, $targetMdDisplay) | Out-Null
$iCellDeDe.Font.Underline = 2
$iCellDeDe.Font.Color = $hyperlinkColor

$wsDeDe.Range("J5").Value = $deDeTargetFile
$wsDeDe.Range("K5").Value = $deDeHandbackDate
$wsDeDe.Range("P5").Value = $errorDetail

$wsDeDe.Columns.Item(16).ColumnWidth = $wideColumnWidth
